$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 18: "indexYear" - percent deviation of population from year 2000 (col AC) ---
$ws.Range("D18").Value = "indexYear"
$ws.Range("I18").Formula = "=(I4 - `$AC`$4)/`$AC`$4"
$ws.Range("J18:AY18").Formula = "=(J4 - `$AC`$4)/`$AC`$4"

# --- Update the saved selection to match the new active cell ---
$ws.Range("J19").Select() | Out-Null
